# Refined metadata to be additional tab
# 1. Update the per-row "time_taken" timestamps on the "data" sheet (F2:F59).
# 2. Add a new "metadata" worksheet after "data" summarising the panel query.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Refresh the time_taken column (F2:F59) on the "data" sheet ---------
$newTimes = @(
    "2021-10-05 14:19:06.461689",
    "2021-10-05 14:19:06.461697",
    "2021-10-05 14:19:06.461700",
    "2021-10-05 14:19:06.461702",
    "2021-10-05 14:19:06.461705",
    "2021-10-05 14:19:06.461708",
    "2021-10-05 14:19:06.461710",
    "2021-10-05 14:19:06.461713",
    "2021-10-05 14:19:06.461716",
    "2021-10-05 14:19:06.461718",
    "2021-10-05 14:19:06.461721",
    "2021-10-05 14:19:06.461723",
    "2021-10-05 14:19:06.461726",
    "2021-10-05 14:19:06.461728",
    "2021-10-05 14:19:06.461731",
    "2021-10-05 14:19:06.461733",
    "2021-10-05 14:19:06.461736",
    "2021-10-05 14:19:06.461738",
    "2021-10-05 14:19:06.461741",
    "2021-10-05 14:19:06.461743",
    "2021-10-05 14:19:06.461746",
    "2021-10-05 14:19:06.461748",
    "2021-10-05 14:19:06.461751",
    "2021-10-05 14:19:06.461753",
    "2021-10-05 14:19:06.461756",
    "2021-10-05 14:19:06.461759",
    "2021-10-05 14:19:06.461761",
    "2021-10-05 14:19:06.461764",
    "2021-10-05 14:19:06.461766",
    "2021-10-05 14:19:06.461769",
    "2021-10-05 14:19:06.461771",
    "2021-10-05 14:19:06.461774",
    "2021-10-05 14:19:06.461777",
    "2021-10-05 14:19:06.461779",
    "2021-10-05 14:19:06.461782",
    "2021-10-05 14:19:06.461784",
    "2021-10-05 14:19:06.461786",
    "2021-10-05 14:19:06.461789",
    "2021-10-05 14:19:06.461791",
    "2021-10-05 14:19:06.461794",
    "2021-10-05 14:19:06.461796",
    "2021-10-05 14:19:06.461799",
    "2021-10-05 14:19:06.461801",
    "2021-10-05 14:19:06.461804",
    "2021-10-05 14:19:06.461806",
    "2021-10-05 14:19:06.461809",
    "2021-10-05 14:19:06.461811",
    "2021-10-05 14:19:06.461814",
    "2021-10-05 14:19:06.461816",
    "2021-10-05 14:19:06.461819",
    "2021-10-05 14:19:06.461821",
    "2021-10-05 14:19:06.461823",
    "2021-10-05 14:19:06.461826",
    "2021-10-05 14:19:06.461829",
    "2021-10-05 14:19:06.461831",
    "2021-10-05 14:19:06.461834",
    "2021-10-05 14:19:06.461836",
    "2021-10-05 14:19:06.461839"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add the "metadata" sheet right after "data" ------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Keep the same outline defaults the rest of the workbook uses.
$meta.Outline.SummaryRow = 1
$meta.Outline.SummaryColumn = 1

# Header row (B1:G1) - reuse the bold/centered header style from "data"!B1:F1.
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row 2 - A2 reuses "data"!A2's style (bordered/centered).
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Adult solid tumours for rare disease"
$meta.Range("C2").Value = 391

# data_version must stay a literal text value ("1.25"), not a number.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.25"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-07-28T13:54:48.174942Z"
$meta.Range("F2").Value = "2021-10-05 14:19:06.457914"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/391/?format=json"

$data.Select()
